# Append a new transaction row (row 36) to the sales log on Sheet1.
# Mirrors the existing "no member / no inventory item" rows (e.g. row 4,
# row 33) where columns B-F are blank text cells and G is the numeric
# total_amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 36

$ws.Cells.Item($row, 1).Value = "2025-03-13 01:07:34"

# Columns B-F (member_id, member_name, inventory_id, item_name, quantity)
# are blank in this transaction, same as the rest of the sheet. A plain
# Value = "" clears the cell entirely instead of leaving a blank text
# cell, so write a lone apostrophe (forces a text-typed empty string,
# matching the workbook's convention) and then strip the resulting
# quote-prefix formatting so no stray style is left on the cell.
2..6 | ForEach-Object {
    $cell = $ws.Cells.Item($row, $_)
    $cell.Value = "'"
    $cell.Style = "Normal"
}

$ws.Cells.Item($row, 7).Value = 75
$ws.Cells.Item($row, 8).Value = "Card"
$ws.Cells.Item($row, 9).Value = "admin"
$ws.Cells.Item($row, 10).Value = "Admin"
$ws.Cells.Item($row, 11).Value = "Badam (1)"
$ws.Cells.Item($row, 12).Value = '[{"name": "Badam", "quantity": 1, "price": 75.0, "total": 75.0}]'
